$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the example row (row 2) - order matters so new shared strings
# land at the same indices as the authored workbook (Example=5, G01 - Rent
# 12/2002=6, G01 - Lease=7, Paid=8).
$ws.Range("F2").Value = "Example"
$ws.Range("A2").Value = "G01 - Rent 12/2002"
$ws.Range("B2").Value = "G01 - Lease"
$ws.Range("C2").Value = 1000
$ws.Range("D2").Value = "11/29/2020"
$ws.Range("E2").Value = "Paid"

# Make row 2 a bit taller to fit the wrapped example text.
$ws.Rows(2).RowHeight = 30

# Highlight row 3 (plus the new column F) in red to flag it in the template.
$ws.Range("A3:F3").Interior.Color = 255

# Move the active selection to A3, matching the saved view state.
[void]$ws.Range("A3").Select()
